$wb = $excel.ActiveWorkbook

$wsDatatypes = $wb.Worksheets.Item("Datatypes")
$wsResults = $wb.Worksheets.Item("SpreadsheetResults")

# --- Datatypes sheet: widen column D and add a new (4th) data row that
# mirrors the existing field row, closing the table off with a bottom border ---
$wsDatatypes.Columns.Item(4).ColumnWidth = 20.66

$wsDatatypes.Range("B4:D4").Copy($wsDatatypes.Range("B5:D5"))
$wsDatatypes.Range("B5:D5").Borders.Item(9).LineStyle = 1
$wsDatatypes.Range("B5:D5").Borders.Item(9).Weight = 2

# --- Selection / active-sheet bookkeeping ---
# First, touch the SpreadsheetResults sheet's selection while it is active ...
$wsResults.Activate()
$wsResults.Range("C6").Select()

# ... then make Datatypes the active (selected) sheet/tab, matching the new
# selection left there.
$wsDatatypes.Activate()
$wsDatatypes.Range("D10").Select()
